$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("J2").Value = ""

# Row 3
$ws.Range("D3").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("G3").Value = "30,46 TL - 60,94 TL - 609,43 TL"

# Row 4
$ws.Range("D4").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("G4").Value = "30,46 TL - 60,94 TL - 609,43 TL"

# Row 5
$ws.Range("D5").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("G5").Value = "30,46 TL - 60,94 TL - 609,43 TL"

# Row 6
$ws.Range("D6").Value = "6,09 TL - 12,19 TL - 152,35 TL"
$ws.Range("G6").Value = "6,09 TL - 12,19 TL - 152,35 TL"

# Row 7
$ws.Range("J7").Value = ""

# Row 8
$ws.Range("D8").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("G8").Value = "15,23 TL - 30,47 TL - 211,05 TL"

# Row 9
$ws.Range("D9").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("G9").Value = "15,23 TL - 30,47 TL - 211,05 TL"

# Row 10
$ws.Range("D10").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("G10").Value = "15,23 TL - 30,47 TL - 211,05 TL"

# Row 11
$ws.Range("D11").Value = "3,04 TL - 6,09 TL - 76,17 TL"
$ws.Range("G11").Value = "3,04 TL - 6,09 TL - 76,17 TL"

# Row 12
$ws.Range("D12").Value = "WU: 0,75 USD–12 USD; Diğer: 700 TL–4.000 TL"
$ws.Range("G12").Value = "Şube (Kasadan): %0,5; Şube (Hesaptan): %0,75; İnternet: 15 USD"

# Row 13
$ws.Range("D13").Value = "Hesaba: Asgari 1 TL | Azami 300 TL"
$ws.Range("E13").Value = "Hesaba: Asgari 1 TL | Azami 851,5 TL"

# Row 14
$ws.Range("D14").Value = "2.300 TL - 9.500 TL"
$ws.Range("G14").Value = "6.300 TL - 6,09 TL"
